# "make sure it works" -- continue the design-notes workbook:
#  - Sheet2's old DB-schema scratch table / "next:" TODO is cleared out
#  - Sheet1 gets the expanded plan (what's shown, then a red TODO list)
#  - Sheet1 becomes the active tab instead of Sheet2

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# Wipe everything that used to live on Sheet2 (DB field names, the
# id/words pairing table, "дальше:" + its single TODO line) while keeping
# the sheet itself / its row grid intact.
$ws2.Range("A3:D36").ClearContents()

# New content continuing the notes on Sheet1.
$ws1.Range("A19").Value = "показывается слово выбранное с шансом из слов пользователя из базы"
$ws1.Range("A20").Value = "показывается на рандомном языке"

# Blank spacer row (kept, just with the default row height set explicitly).
$ws1.Rows.Item(21).RowHeight = 14.25

$ws1.Range("A22").Value = "ниже показываются клетки для ввода (кол-во букв и пробелов в ответе известно)"

# The "дальше:" follow-up items, re-done as a highlighted (red) TODO list.
$ws1.Range("A25").Value = "добавить проверку на правильный/неправильный ответ + реакцию БД"
$ws1.Range("A26").Value = "добавить отображение статистики (слов изучено по категориям / всего)"
$ws1.Range("A27").Value = "добавить процент изучения слова которое выпало"
$ws1.Range("A28").Value = "добавить кнопки для того чтобы отмечать категории которые нужны"
$ws1.Range("A29").Value = "сделать базовый html шаблон и к нему прикрутить остальные"
$ws1.Range("A30").Value = "сделать навигацию на главной странице"
$ws1.Range("A31").Value = "сделать страницу админку со статистикой"
$ws1.Range("A32").Value = "сделать скрипт для загрузки новых слов из эксель файла"

# ColorIndex 3 == the OOXML legacy-palette "indexed 2" red (FF0000); the
# Excel ColorIndex palette is 1-based / offset by one from OOXML's.
$ws1.Range("A25:A32").Font.ColorIndex = 3

# Sheet1 is now the active tab (was Sheet2).
$ws1.Activate()
